# Auto-generated edit script applying the Mateus_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 177.42857
$ws.Range("I2").Value = 92.96
$ws.Range("J2").Value = 881.3333
$ws.Range("K2").Value = 92.96
$ws.Range("L2").Value = 881.3333
$ws.Range("M2").Value = 20.04000000000001
$ws.Range("N2").Value = -1107.3333
$ws.Range("H17").Value = 3496.9092
$ws.Range("J17").Value = 3496.9092
$ws.Range("L17").Value = 10490.7276
$ws.Range("N17").Value = -10826.7276
$ws.Range("H19").Value = 3698.1
$ws.Range("I19").Value = 1328
$ws.Range("K19").Value = 1328
$ws.Range("M19").Value = -1153
$ws.Range("H34").Value = 5999.6665
$ws.Range("I34").Value = 5999.6665
$ws.Range("K34").Value = 5999.6665
$ws.Range("M34").Value = -5796.6665
$ws.Range("H36").Value = 5999.6665
$ws.Range("I36").Value = 5999.6665
$ws.Range("K36").Value = 5999.6665
$ws.Range("M36").Value = -5284.6665
$ws.Range("H40").Value = 16672458
$ws.Range("J40").Value = 22733808
$ws.Range("L40").Value = 22733808
$ws.Range("N40").Value = -22734158
$ws.Range("H55").Value = 577.6923
$ws.Range("J55").Value = 803.5
$ws.Range("L55").Value = 803.5
$ws.Range("N55").Value = -1231.5
$ws.Range("H62").Value = 5930.1
$ws.Range("I62").Value = 5977.8887
$ws.Range("K62").Value = 5977.8887
$ws.Range("M62").Value = -5353.8887
$ws.Range("H64").Value = 6687.278
$ws.Range("I64").Value = 4976.8335
$ws.Range("J64").Value = 10108.167
$ws.Range("K64").Value = 4976.8335
$ws.Range("L64").Value = 10108.167
$ws.Range("M64").Value = -4728.8335
$ws.Range("N64").Value = -10604.167
$ws.Range("H65").Value = 5930.1
$ws.Range("I65").Value = 5977.8887
$ws.Range("K65").Value = 29889.4435
$ws.Range("M65").Value = -26769.4435
$ws.Range("H67").Value = 6687.278
$ws.Range("I67").Value = 4976.8335
$ws.Range("J67").Value = 10108.167
$ws.Range("K67").Value = 4976.8335
$ws.Range("L67").Value = 10108.167
$ws.Range("M67").Value = -4118.8335
$ws.Range("N67").Value = -11824.167
$ws.Range("H70").Value = 3391.4285
$ws.Range("I70").Value = 2185.5
$ws.Range("K70").Value = 6556.5
$ws.Range("M70").Value = -6286.5
$ws.Range("H73").Value = 3391.4285
$ws.Range("I73").Value = 2185.5
$ws.Range("K73").Value = 6556.5
$ws.Range("M73").Value = -5620.5
$ws.Range("H86").Value = 2650
$ws.Range("I86").Value = 2650
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2650
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1527
$ws.Range("N86").ClearContents()
$ws.Range("H88").Value = 1713
$ws.Range("I88").Value = 1679.1428
$ws.Range("J88").Value = 1739.3334
$ws.Range("K88").Value = 1679.1428
$ws.Range("L88").Value = 1739.3334
$ws.Range("M88").Value = -1273.1428
$ws.Range("N88").Value = -2551.3334
$ws.Range("H89").Value = 2650
$ws.Range("I89").Value = 2650
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 13250
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -7634
$ws.Range("N89").ClearContents()
$ws.Range("H91").Value = 1713
$ws.Range("I91").Value = 1679.1428
$ws.Range("J91").Value = 1739.3334
$ws.Range("K91").Value = 1679.1428
$ws.Range("L91").Value = 1739.3334
$ws.Range("M91").Value = -275.1428000000001
$ws.Range("N91").Value = -4547.3334
$ws.Range("H98").Value = 1504.0476
$ws.Range("I98").Value = 1504.0476
$ws.Range("K98").Value = 1504.0476
$ws.Range("M98").Value = -6.047600000000102
$ws.Range("H112").Value = 2284.0833
$ws.Range("J112").Value = 2318.0908
$ws.Range("L112").Value = 6954.2724
$ws.Range("N112").Value = -9170.2724
$ws.Range("H122").Value = 1504.0476
$ws.Range("I122").Value = 1504.0476
$ws.Range("K122").Value = 4512.142800000001
$ws.Range("M122").Value = -2062.142800000001
$ws.Range("H132").Value = 881.2553
$ws.Range("I132").Value = 904.5333
$ws.Range("K132").Value = 2713.5999
$ws.Range("M132").Value = -183.5999000000002
$ws.Range("H135").Value = 1423.6666
$ws.Range("I135").Value = 1453.9286
$ws.Range("K135").Value = 13085.3574
$ws.Range("M135").Value = -10550.3574
$ws.Range("H137").Value = 6450.8887
$ws.Range("I137").Value = 9122.9
$ws.Range("J137").Value = 3110.875
$ws.Range("K137").Value = 27368.7
$ws.Range("L137").Value = 9332.625
$ws.Range("M137").Value = -24818.7
$ws.Range("N137").Value = -14432.625
$ws.Range("H138").Value = 5374.1084
$ws.Range("I138").Value = 6635.857
$ws.Range("J138").Value = 4731.7637
$ws.Range("K138").Value = 19907.571
$ws.Range("L138").Value = 14195.2911
$ws.Range("M138").Value = -14767.571
$ws.Range("N138").Value = -24475.2911
$ws.Range("H141").Value = 1934.5264
$ws.Range("J141").Value = 2150
$ws.Range("L141").Value = 6450
$ws.Range("N141").Value = -16810

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26024.857
$ws.Range("I32").Value = 26024.857
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 26024.857
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -25737.857
$ws.Range("N32").ClearContents()
$ws.Range("H45").Value = 3508.889
$ws.Range("I45").Value = 3058.818
$ws.Range("J45").Value = 5489.2
$ws.Range("K45").Value = 3058.818
$ws.Range("L45").Value = 5489.2
$ws.Range("M45").Value = -2681.818
$ws.Range("N45").Value = -6243.2
$ws.Range("H61").Value = 6927.9653
$ws.Range("I61").Value = 6810.7085
$ws.Range("J61").Value = 7490.8
$ws.Range("K61").Value = 6810.7085
$ws.Range("L61").Value = 7490.8
$ws.Range("M61").Value = -6598.7085
$ws.Range("N61").Value = -7914.8
$ws.Range("H103").Value = 39998.5
$ws.Range("J103").Value = 39998.5
$ws.Range("L103").Value = 39998.5
$ws.Range("N103").Value = -42342.5
$ws.Range("H110").Value = 3641.3845
$ws.Range("I110").Value = 1939.8182
$ws.Range("K110").Value = 1939.8182
$ws.Range("M110").Value = 105.1818000000001
$ws.Range("H122").Value = 2502.5557
$ws.Range("I122").Value = 2646.2856
$ws.Range("J122").Value = 1999.5
$ws.Range("K122").Value = 7938.8568
$ws.Range("L122").Value = 5998.5
$ws.Range("M122").Value = -5488.8568
$ws.Range("N122").Value = -10898.5
$ws.Range("H132").Value = 3078.0833
$ws.Range("I132").Value = 2942
$ws.Range("K132").Value = 8826
$ws.Range("M132").Value = -6296
$ws.Range("H136").Value = 6927.9653
$ws.Range("I136").Value = 6810.7085
$ws.Range("J136").Value = 7490.8
$ws.Range("K136").Value = 20432.1255
$ws.Range("L136").Value = 22472.4
$ws.Range("M136").Value = -17882.1255
$ws.Range("N136").Value = -27572.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3165.6
$ws.Range("I20").Value = 2753.9
$ws.Range("J20").Value = 3989
$ws.Range("K20").Value = 2753.9
$ws.Range("L20").Value = 3989
$ws.Range("M20").Value = -2506.9
$ws.Range("N20").Value = -4483
$ws.Range("H94").Value = 1178.5714
$ws.Range("I94").Value = 1209.25
$ws.Range("K94").Value = 1209.25
$ws.Range("M94").Value = -758.25
$ws.Range("H97").Value = 13652.833
$ws.Range("I97").Value = 13652.833
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 13652.833
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -12661.833
$ws.Range("N97").ClearContents()
$ws.Range("H105").Value = 1978.3125
$ws.Range("I105").Value = 1959
$ws.Range("K105").Value = 1959
$ws.Range("M105").Value = -212
$ws.Range("H107").Value = 86298.5
$ws.Range("I107").Value = 2676.4443
$ws.Range("J107").Value = 337164.66
$ws.Range("K107").Value = 2676.4443
$ws.Range("L107").Value = 337164.66
$ws.Range("M107").Value = -756.4443000000001
$ws.Range("N107").Value = -341004.66
$ws.Range("H134").Value = 4844.8965
$ws.Range("I134").Value = 4875.0713
$ws.Range("K134").Value = 14625.2139
$ws.Range("M134").Value = -12090.2139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2061.3684
$ws.Range("I16").Value = 1929.25
$ws.Range("J16").Value = 2287.8572
$ws.Range("K16").Value = 1929.25
$ws.Range("L16").Value = 2287.8572
$ws.Range("M16").Value = -1642.25
$ws.Range("N16").Value = -2861.8572
$ws.Range("H19").Value = 190
$ws.Range("I19").Value = 170.36363
$ws.Range("J19").Value = 406
$ws.Range("K19").Value = 170.36363
$ws.Range("L19").Value = 406
$ws.Range("M19").Value = -0.3636300000000006
$ws.Range("N19").Value = -746
$ws.Range("H24").Value = 190
$ws.Range("I24").Value = 170.36363
$ws.Range("J24").Value = 406
$ws.Range("K24").Value = 170.36363
$ws.Range("L24").Value = 406
$ws.Range("M24").Value = -0.3636300000000006
$ws.Range("N24").Value = -746
$ws.Range("H31").Value = 5670.5
$ws.Range("I31").Value = 3238.4443
$ws.Range("J31").Value = 12966.667
$ws.Range("K31").Value = 3238.4443
$ws.Range("L31").Value = 12966.667
$ws.Range("M31").Value = -2943.4443
$ws.Range("N31").Value = -13556.667
$ws.Range("H34").Value = 5670.5
$ws.Range("I34").Value = 3238.4443
$ws.Range("J34").Value = 12966.667
$ws.Range("K34").Value = 3238.4443
$ws.Range("L34").Value = 12966.667
$ws.Range("M34").Value = -3036.4443
$ws.Range("N34").Value = -13370.667
$ws.Range("H41").Value = 29849.5
$ws.Range("I41").Value = 35700
$ws.Range("K41").Value = 35700
$ws.Range("M41").Value = -35272
$ws.Range("H43").Value = 22222
$ws.Range("J43").Value = 22222
$ws.Range("L43").Value = 22222
$ws.Range("N43").Value = -22590
$ws.Range("H50").Value = 10000
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 15122.25
$ws.Range("H58").Value = 4592.5
$ws.Range("I58").Value = 3044.125
$ws.Range("J58").Value = 7069.9
$ws.Range("K58").Value = 3044.125
$ws.Range("L58").Value = 7069.9
$ws.Range("M58").Value = -2841.125
$ws.Range("N58").Value = -7475.9
$ws.Range("H60").Value = 40000
$ws.Range("H61").Value = 15122.25
$ws.Range("H99").Value = 5162.647
$ws.Range("J99").Value = 5942.364
$ws.Range("L99").Value = 5942.364
$ws.Range("N99").Value = -8938.364
$ws.Range("H101").Value = 22222
$ws.Range("J101").Value = 22222
$ws.Range("L101").Value = 22222
$ws.Range("N101").Value = -28712
$ws.Range("H105").Value = 1418.2051
$ws.Range("I105").Value = 1064.1945
$ws.Range("K105").Value = 1064.1945
$ws.Range("M105").Value = 682.8054999999999
$ws.Range("H107").Value = 116099.22
$ws.Range("I107").Value = 1981.6
$ws.Range("K107").Value = 1981.6
$ws.Range("M107").Value = -61.59999999999991
$ws.Range("H113").Value = 2061.3684
$ws.Range("I113").Value = 1929.25
$ws.Range("J113").Value = 2287.8572
$ws.Range("K113").Value = 1929.25
$ws.Range("L113").Value = 2287.8572
$ws.Range("M113").Value = 240.75
$ws.Range("N113").Value = -6627.8572
$ws.Range("H126").Value = 5162.647
$ws.Range("J126").Value = 5942.364
$ws.Range("L126").Value = 17827.092
$ws.Range("N126").Value = -22767.092
$ws.Range("H132").Value = 5230.472
$ws.Range("I132").Value = 5746.04
$ws.Range("K132").Value = 17238.12
$ws.Range("M132").Value = -14708.12
$ws.Range("H134").Value = 2791.3057
$ws.Range("I134").Value = 1681.4546
$ws.Range("K134").Value = 5044.3638
$ws.Range("M134").Value = -2509.3638
$ws.Range("H136").Value = 4592.5
$ws.Range("I136").Value = 3044.125
$ws.Range("J136").Value = 7069.9
$ws.Range("K136").Value = 9132.375
$ws.Range("L136").Value = 21209.7
$ws.Range("M136").Value = -6582.375
$ws.Range("N136").Value = -26309.7
$ws.Range("H141").Value = 259033.3
$ws.Range("J141").Value = 310416.7
$ws.Range("L141").Value = 310416.7
$ws.Range("N141").Value = -320776.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 407.5
$ws.Range("J34").Value = 750
$ws.Range("L34").Value = 2250
$ws.Range("N34").Value = -2418
$ws.Range("H39").Value = 7313.3335
$ws.Range("J39").Value = 8141.6665
$ws.Range("L39").Value = 24424.9995
$ws.Range("N39").Value = -25012.9995
$ws.Range("H45").Value = 2766.5
$ws.Range("I45").Value = 999.5
$ws.Range("K45").Value = 2998.5
$ws.Range("M45").Value = -2466.5
$ws.Range("H55").Value = 14525
$ws.Range("J55").Value = 14525
$ws.Range("L55").Value = 43575
$ws.Range("N55").Value = -43929
$ws.Range("H86").Value = 296
$ws.Range("I86").Value = 296
$ws.Range("K86").Value = 888
$ws.Range("M86").Value = 298
$ws.Range("H89").Value = 296
$ws.Range("I89").Value = 296
$ws.Range("K89").Value = 2664
$ws.Range("M89").Value = 3264
$ws.Range("H92").Value = 780.6667
$ws.Range("I92").Value = 177.33333
$ws.Range("J92").Value = 1082.3334
$ws.Range("K92").Value = 531.99999
$ws.Range("L92").Value = 3247.0002
$ws.Range("M92").Value = 716.00001
$ws.Range("N92").Value = -5743.0002
$ws.Range("H103").Value = 410.16666
$ws.Range("I103").Value = 410.16666
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 1230.49998
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -351.4999800000001
$ws.Range("N103").ClearContents()
$ws.Range("H107").Value = 1289
$ws.Range("I107").Value = 623.75
$ws.Range("J107").Value = 3950
$ws.Range("K107").Value = 1871.25
$ws.Range("L107").Value = 11850
$ws.Range("M107").Value = 48.75
$ws.Range("N107").Value = -15690
$ws.Range("H114").Value = 4999.5
$ws.Range("J114").Value = 5000
$ws.Range("L114").Value = 15000
$ws.Range("N114").Value = -21508
$ws.Range("H122").Value = 18037.05
$ws.Range("I122").Value = 8166.6665
$ws.Range("K122").Value = 73499.9985
$ws.Range("M122").Value = -71049.9985
$ws.Range("H128").Value = 219869.88
$ws.Range("I128").Value = 219869.88
$ws.Range("K128").Value = 659609.64
$ws.Range("M128").Value = -654629.64

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 89.82353
$ws.Range("I2").Value = 101.6
$ws.Range("J2").Value = 1.5
$ws.Range("K2").Value = 101.6
$ws.Range("L2").Value = 1.5
$ws.Range("M2").Value = 11.40000000000001
$ws.Range("N2").Value = -227.5
$ws.Range("H70").Value = 9198.6
$ws.Range("I70").Value = 7372.5
$ws.Range("J70").Value = 11285.571
$ws.Range("K70").Value = 7372.5
$ws.Range("L70").Value = 11285.571
$ws.Range("M70").Value = -7102.5
$ws.Range("N70").Value = -11825.571
$ws.Range("H73").Value = 9198.6
$ws.Range("I73").Value = 7372.5
$ws.Range("J73").Value = 11285.571
$ws.Range("K73").Value = 7372.5
$ws.Range("L73").Value = 11285.571
$ws.Range("M73").Value = -6436.5
$ws.Range("N73").Value = -13157.571
$ws.Range("H80").Value = 2319.9
$ws.Range("I80").Value = 2396.8333
$ws.Range("K80").Value = 2396.8333
$ws.Range("M80").Value = -1398.8333
$ws.Range("H83").Value = 2319.9
$ws.Range("I83").Value = 2396.8333
$ws.Range("K83").Value = 11984.1665
$ws.Range("M83").Value = -6992.166499999999
$ws.Range("H97").Value = 918.4286
$ws.Range("I97").Value = 903.4545
$ws.Range("J97").Value = 973.3333
$ws.Range("K97").Value = 903.4545
$ws.Range("L97").Value = 973.3333
$ws.Range("M97").Value = -407.4545000000001
$ws.Range("N97").Value = -1965.3333
$ws.Range("H122").Value = 5264.706
$ws.Range("I122").Value = 5057.6665
$ws.Range("J122").Value = 6817.5
$ws.Range("K122").Value = 15172.9995
$ws.Range("L122").Value = 20452.5
$ws.Range("M122").Value = -12722.9995
$ws.Range("N122").Value = -25352.5
$ws.Range("H132").Value = 2562.4827
$ws.Range("I132").Value = 2562.4827
$ws.Range("K132").Value = 7687.4481
$ws.Range("M132").Value = -5157.4481
$ws.Range("H135").Value = 149165.83
$ws.Range("J135").Value = 149165.83
$ws.Range("L135").Value = 149165.83
$ws.Range("N135").Value = -159305.83

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 2574.8333
$ws.Range("J19").Value = 3399.5
$ws.Range("L19").Value = 3399.5
$ws.Range("N19").Value = -3739.5
$ws.Range("H22").Value = 6196.1113
$ws.Range("J22").Value = 14996.333
$ws.Range("L22").Value = 14996.333
$ws.Range("N22").Value = -15586.333
$ws.Range("H27").Value = 6196.1113
$ws.Range("J27").Value = 14996.333
$ws.Range("L27").Value = 14996.333
$ws.Range("N27").Value = -15210.333
$ws.Range("H40").Value = 6544.9
$ws.Range("I40").Value = 6978.4287
$ws.Range("K40").Value = 6978.4287
$ws.Range("M40").Value = -6842.4287
$ws.Range("H46").Value = 4608.2856
$ws.Range("I46").Value = 2924.1
$ws.Range("J46").Value = 8818.75
$ws.Range("K46").Value = 2924.1
$ws.Range("L46").Value = 8818.75
$ws.Range("M46").Value = -2736.1
$ws.Range("N46").Value = -9194.75
$ws.Range("H68").Value = 6499.5
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 6499.5
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H100").Value = 11606.2
$ws.Range("I100").Value = 4624
$ws.Range("J100").Value = 16261
$ws.Range("K100").Value = 4624
$ws.Range("L100").Value = 16261
$ws.Range("M100").Value = -4083
$ws.Range("N100").Value = -17343
$ws.Range("H106").Value = 18666.334
$ws.Range("J106").Value = 18666.334
$ws.Range("L106").Value = 18666.334
$ws.Range("N106").Value = -21190.334
$ws.Range("H122").Value = 6609.6875
$ws.Range("I122").Value = 6609.6875
$ws.Range("K122").Value = 19829.0625
$ws.Range("M122").Value = -17379.0625
$ws.Range("H132").Value = 13507.227
$ws.Range("I132").Value = 12860.934
$ws.Range("K132").Value = 38582.802
$ws.Range("M132").Value = -36052.802

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 12000
$ws.Range("J21").Value = 12000
$ws.Range("L21").Value = 12000
$ws.Range("N21").Value = -12470
$ws.Range("H29").Value = 29.5
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H35").Value = 12000
$ws.Range("J35").Value = 12000
$ws.Range("L35").Value = 12000
$ws.Range("N35").Value = -12580
$ws.Range("H46").Value = 78304.625
$ws.Range("J46").Value = 78304.625
$ws.Range("L46").Value = 78304.625
$ws.Range("N46").Value = -78766.625
$ws.Range("H54").Value = 133500
$ws.Range("J54").Value = 250000
$ws.Range("L54").Value = 250000
$ws.Range("N54").Value = -251040
$ws.Range("H81").Value = 4204.5835
$ws.Range("I81").Value = 4135.1816
$ws.Range("K81").Value = 8270.3632
$ws.Range("M81").Value = -7209.3632
$ws.Range("H84").Value = 4204.5835
$ws.Range("I84").Value = 4135.1816
$ws.Range("K84").Value = 41351.816
$ws.Range("M84").Value = -36047.816
$ws.Range("H100").Value = 1375.6666
$ws.Range("I100").Value = 1000
$ws.Range("J100").Value = 1422.625
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 2845.25
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -3927.25
$ws.Range("H107").Value = 25982.195
$ws.Range("I107").Value = 1233.56
$ws.Range("K107").Value = 3700.68
$ws.Range("M107").Value = -1780.68
$ws.Range("H113").Value = 4156.6665
$ws.Range("I113").Value = 1719.3846
$ws.Range("K113").Value = 5158.1538
$ws.Range("M113").Value = -2988.1538
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800
$ws.Range("H124").Value = 93525.836
$ws.Range("J124").Value = 93525.836
$ws.Range("L124").Value = 93525.836
$ws.Range("N124").Value = -103345.836
$ws.Range("H126").Value = 3720.375
$ws.Range("I126").Value = 2691.2104
$ws.Range("K126").Value = 8073.6312
$ws.Range("M126").Value = -5603.6312
$ws.Range("H132").Value = 4412.3726
$ws.Range("I132").Value = 3716.175
$ws.Range("J132").Value = 6944
$ws.Range("K132").Value = 11148.525
$ws.Range("L132").Value = 20832
$ws.Range("M132").Value = -8618.525000000001
$ws.Range("N132").Value = -25892
$ws.Range("H134").Value = 78304.625
$ws.Range("J134").Value = 78304.625
$ws.Range("L134").Value = 234913.875
$ws.Range("N134").Value = -239983.875
$ws.Range("H136").Value = 3208.6052
$ws.Range("I136").Value = 2151.1333
$ws.Range("K136").Value = 6453.3999
$ws.Range("M136").Value = -3903.3999

Write-Host "Applied Mateus_Profits changes"